$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Datos actualizados" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 04:34"

# --- 2. Read the full data block (countries table, rows 4-218, cols A-H) ---
$firstRow = 4
$lastRow = 218
$rng = $ws.Range("A$firstRow`:H$lastRow")
$data = $rng.Value()
$rowCount = $data.GetLength(0)

# --- 3. Update individual country stats that changed ---
for ($r = 1; $r -le $rowCount; $r++) {
    $country = $data[$r, 1]

    if ($country -eq "Guatemala") {
        $data[$r, 2] = 900
        $data[$r, 3] = 68
        $data[$r, 4] = 101
        $data[$r, 5] = 775
        $data[$r, 6] = 5
        $data[$r, 7] = 1
        $data[$r, 8] = 24
    }
    elseif ($country -eq "Guayana Francesa") {
        $data[$r, 4] = 119
        $data[$r, 5] = 21
    }
}

# --- 4. Re-sort the table by "Casos totales" (column B) descending, stable ---
$rowsList = New-Object System.Collections.ArrayList
for ($r = 1; $r -le $rowCount; $r++) {
    $rowVals = New-Object System.Collections.ArrayList
    for ($c = 1; $c -le 8; $c++) {
        [void]$rowVals.Add($data[$r, $c])
    }
    [void]$rowsList.Add([pscustomobject]@{ Order = $r; Cols = $rowVals; Key = [double]$data[$r, 2] })
}

$sortedRows = $rowsList | Sort-Object -Property Key -Descending

# New-Object 'object[,]' yields a 0-based .NET array (unlike the 1-based
# array returned by Range.Value()), so index it from 0.
$outData = New-Object 'object[,]' $rowCount,8
$r = 0
foreach ($item in $sortedRows) {
    for ($c = 0; $c -le 7; $c++) {
        $outData[$r, $c] = $item.Cols[$c]
    }
    $r++
}

$rng.Value = $outData
